# Updated cryptos list on Fri Aug 18 21:36:58 UTC 2023 with GitHub Actions
# Refreshes Price (D) and Volume(1h) (E) values for each coin row, and
# reflects that the Polkadot / WrappedEther rows (12 and 13) swapped order.
# NumberFormat is forced to Text ("@") before assigning, because several
# values (e.g. "218.30") would otherwise be auto-parsed as numbers by Excel,
# which strips significant trailing/duplicate-decimal formatting that the
# source data relies on (dotted thousands separators, trailing zeros, etc).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.257.58"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -5.42%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.674.15"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -2.73%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.45%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.30"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5116"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -8.95%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.37%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2661"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06396"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.91%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.53"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -4.65%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.59%  "
$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.570"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.41%  "
$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.674.12"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -2.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5831"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.14%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.903.95"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.58%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008688"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.28%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.87"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -11.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.354.78"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -4.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.966"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.96%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.17%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "190.64"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.11%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.235"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.79%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.007"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.28"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.21%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.679"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1183"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.92%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.69"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.41%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05908"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.41%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.279"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -6.10%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.324"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -3.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.534"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -3.21%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.522"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -3.60%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.648"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.81%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.014"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6031"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -5.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.362"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.45%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.653"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.84%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.047"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.080.46"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -2.22%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8695"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.43%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.83"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.825.04"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000113"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.10"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.06%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.36%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.098"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4300"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05194"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.97%  "
